$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# -----------------------------------------------------------------
# 1. Insert a blank column at D so "urban_pop" has a home right after
#    "rural_urban" (C). This shifts D..Y -> E..Z and slides the
#    existing merged header groups (D1:I1 etc.) one column right.
# -----------------------------------------------------------------
$ws.Columns.Item(4).Insert()

# -----------------------------------------------------------------
# 2. Insert a second blank column right after "multiplier" (now M)
#    to hold the new "hivalu" field, before "dep" (now N).
# -----------------------------------------------------------------
$ws.Columns.Item(14).Insert()

# -----------------------------------------------------------------
# 3. Move "urban_pop" column's data+style from its old position
#    (now column Z, after both inserts) into the new column D.
# -----------------------------------------------------------------
$ws.Range("Z1:Z3").Cut($ws.Range("D1:D3"))

# -----------------------------------------------------------------
# 4. Populate the new "hivalu" header cell (N2) and data cell (N3).
# -----------------------------------------------------------------
$ws.Range("N2").Value = "hivalu"
$ws.Range("N2").Font.Bold = $false
$ws.Range("N3").Formula = ""

# -----------------------------------------------------------------
# 5. Move "Dispersion" label from C1 to B1 so it can anchor a new
#    B1:D1 merge spanning FacilityID/met_station/rural_urban/urban_pop.
# -----------------------------------------------------------------
$ws.Range("B1").Value = $ws.Range("C1").Value()
$ws.Range("C1").Value = ""

# -----------------------------------------------------------------
# 6. Re-box the three "plain" (non-filled) header groups with the
#    composite thin border (left/middle/right split) that Excel
#    produces when you box an unmerged multi-cell range, then merge.
# -----------------------------------------------------------------
$g1 = $ws.Range("B1:D1")
$g1.Borders.LineStyle = 0
$g1.BorderAround(1, 2)
$g1.Merge()

$g2 = $ws.Range("K1:N1")
$g2.UnMerge()
$g2.Borders.LineStyle = 0
$g2.BorderAround(1, 2)
$g2.Merge()

$g3 = $ws.Range("V1:Z1")
$g3.UnMerge()
$g3.Borders.LineStyle = 0
$g3.Interior.Pattern = 0
$g3.BorderAround(1, 2)
$g3.Merge()

# -----------------------------------------------------------------
# 7. Sheet view: drop the frozen/scrolled topLeftCell, move the
#    active selection to V7.
# -----------------------------------------------------------------
$ws.Range("V7").Select()
$excel.ActiveWindow.ScrollColumn = 1
$excel.ActiveWindow.ScrollRow = 1

# -----------------------------------------------------------------
# 8. Page setup: portrait orientation.
# -----------------------------------------------------------------
$ws.PageSetup.Orientation = 1

Write-Host "edit complete"
